# Update "Wind SA Age" table:
#  - rename header "Station" -> "Generator" (A1) (B1 "Commission Year" unchanged text)
#  - insert six additional wind farms (commissioned 2009) in alphabetical order,
#    pushing the existing rows down, and update a handful of commission years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data, in row order, after the edit (header + 19 data rows -> A1:B20)
$data = @(
    @("Generator", "Commission Year"),
    @("Hallett Wind Farm", 2007),
    @("Lake Bonney Stage 2", 2007),
    @("Snowtown Wind Farm", 2008),
    @("Canunda Wind Farm", 2009),
    @("Cathedral Rocks Wind Farm", 2009),
    @("Clements Gap Wind Farm", 2009),
    @("Lake Bonney Stage 1", 2009),
    @("Mount Millar Wind Farm", 2009),
    @("Starfish Hill Wind Farm", 2009),
    @("Wattle Point Wind Farm", 2009),
    @("Lake Bonney Stage 3", 2010),
    @("North Brown Hill Wind Farm", 2010),
    @("Waterloo Wind Farm", 2010),
    @("The Bluff Wind Farm", 2011),
    @("Snowtown South Wind Farm", 2013),
    @("Snowtown North Wind Farm", 2014),
    @("Hornsdale Wind Farm", 2016),
    @("Willogoleche Wind Farm", 2018),
    @("Lincoln Gap Wind Farm", 2019)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A1:B20").Select()
